# A new daily price record was inserted into the weekly series at row 8
# (Terminal Hortofrutícola Agro Chillán - Cilantro). This pushes every
# existing record from row 8 onward down by one row (so old row 8 becomes
# row 9, old row 9 becomes row 10, ..., old row 37 becomes row 38), and the
# new row 8 is populated with a fresh observation dated 2022-04-08 (serial
# 44659) carrying the same Volumen/Precio/Unidad/Origen values as the
# record that used to occupy row 8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 8, shifting rows 8:37 down to 9:38.
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new observation.
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(8, 3).Value = "Ñuble"
$ws.Cells.Item(8, 4).Value = 44659
$ws.Cells.Item(8, 5).Value = 16
$ws.Cells.Item(8, 6).Value = 100112040
$ws.Cells.Item(8, 7).Value = "Cilantro"
$ws.Cells.Item(8, 8).Value = "Sin especificar"
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 200
$ws.Cells.Item(8, 11).Value = 550
$ws.Cells.Item(8, 12).Value = 600
$ws.Cells.Item(8, 13).Value = 575
$ws.Cells.Item(8, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(8, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(8, 16).Value = 575
$ws.Cells.Item(8, 17).Value = 1
$ws.Cells.Item(8, 18).Value = "Hortaliza"
